$d = $word.ActiveDocument

# 1. "log_workout" -> "workout"
$d.Content.Find.Execute("log_workout", $true, $false, $false, $false, $false,
                         $true, 1, $false, "workout", 2)

# 2. Add a new bullet "Borrar registros" right after the existing
#    "Ver registros" bullet (same list style/level as its siblings).
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Ver registros*") {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -gt 0) {
    $anchorPara = $d.Paragraphs($anchorIndex)
    $anchorPara.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs($anchorIndex + 1)
    $newPara.Range.Text = "Borrar registros"
}

# 3. Mark the "Fuentedeprrafopredeter" (Default Paragraph Font) style as
#    semi-hidden ("Hide until used" in Word's Manage Styles dialog), which
#    persists as <w:semiHidden/> on the style definition.
$style = $d.Styles("Fuentedeprrafopredeter")
$style.Hidden = $true
